$d = $word.ActiveDocument

# --- Paragraph 2: "This is the document..." -> "To pull the request into the mainline:"
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "To pull the request into the mainline:"

# --- Paragraph 3: the multi-run "Add instructions..." paragraph becomes the first of
# several new bullet-style paragraphs. Clear its content (keep the paragraph mark).
$p3 = $d.Paragraphs(3)
$p3start = $p3.Range.Start
$p3endNoMark = $p3.Range.End - 1
$clear = $d.Range($p3start, $p3endNoMark)
$clear.Text = ""

# Build paragraph 3 runs: "-Go to repository " + "CMP73010-Ass1-2017"
$ins = $d.Range($p3start, $p3start)
$ins.InsertAfter("-Go to repository ")
$p3 = $d.Paragraphs(3)
$end = $p3.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("CMP73010-Ass1-2017")

# Append a new paragraph mark, then fill paragraph 4 runs:
# "-" + "Click on the pull requests under the repository name" + "."
$p3 = $d.Paragraphs(3)
$end = $p3.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$end = $p4.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("-")
$p4 = $d.Paragraphs(4)
$end = $p4.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("Click on the pull requests under the repository name")
$p4 = $d.Paragraphs(4)
$end = $p4.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter(".")

# New paragraph 5: "-Find pull requested under my username."
$p4 = $d.Paragraphs(4)
$end = $p4.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$end = $p5.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("-Find pull requested under my username.")

# New paragraph 6: "-Click merge pull request."
$p5 = $d.Paragraphs(5)
$end = $p5.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$end = $p6.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("-Click merge pull request.")

# New paragraph 7: "-Enter commit message and click confirm." + bookmark _GoBack
$p6 = $d.Paragraphs(6)
$end = $p6.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertParagraphAfter()
$p7 = $d.Paragraphs(7)
$end = $p7.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertAfter("-Enter commit message and click confirm.")

$p7 = $d.Paragraphs(7)
$end = $p7.Range.End - 1
$bm = $d.Range($end, $end)
$d.Bookmarks.Add("_GoBack", $bm)

# --- Trailing empty paragraph
$p7 = $d.Paragraphs(7)
$end = $p7.Range.End - 1
$ins = $d.Range($end, $end)
$ins.InsertParagraphAfter()
